$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.779.90'
$ws.Range('E2').Value = '  +1.20%  '
$ws.Range('D3').Value = '3.113.15'
$ws.Range('E3').Value = '  -0.10%  '
$ws.Range('E4').Value = '  -0.08%  '
$rng = $ws.Range('D5')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '585.47'
$rng.Style = $origStyle
$ws.Range('E5').Value = '  -0.32%  '
$rng = $ws.Range('D6')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '146.23'
$rng.Style = $origStyle
$ws.Range('E6').Value = '  +1.67%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.107.56'
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -0.47%  '
$ws.Range('E10').Value = '  +10.33%  '
$ws.Range('E11').Value = '  -0.10%  '
$ws.Range('E12').Value = '  -0.75%  '
$ws.Range('E13').Value = '  +2.69%  '
$rng = $ws.Range('D14')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '37.06'
$rng.Style = $origStyle
$ws.Range('E14').Value = '  +4.15%  '
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = '3.631.68'
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').Value = '63.709.46'
$ws.Range('E17').Value = '  +1.14%  '
$rng = $ws.Range('D18')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '7.12'
$rng.Style = $origStyle
$ws.Range('E18').Value = '  -2.16%  '
$ws.Range('D19').Value = '3.116.17'
$ws.Range('E19').Value = '  +0.12%  '
$rng = $ws.Range('D20')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '464.75'
$rng.Style = $origStyle
$ws.Range('E20').Value = '  +2.40%  '
$rng = $ws.Range('D21')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '14.30'
$rng.Style = $origStyle
$ws.Range('E21').Value = '  +1.50%  '
$ws.Range('E22').Value = '  -0.61%  '
$rng = $ws.Range('D23')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '7.52'
$rng.Style = $origStyle
$ws.Range('E23').Value = '  +0.13%  '
$rng = $ws.Range('D24')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '13.09'
$rng.Style = $origStyle
$ws.Range('E24').Value = '  -4.02%  '
$rng = $ws.Range('D25')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '81.99'
$rng.Style = $origStyle
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E27').Value = '  +8.00%  '
$ws.Range('E28').Value = '  -0.66%  '
$ws.Range('E29').Value = '  -1.59%  '
$ws.Range('E30').Value = '  -0.08%  '
$rng = $ws.Range('D31')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '6.85'
$rng.Style = $origStyle
$ws.Range('E31').Value = '  -0.20%  '
$rng = $ws.Range('D32')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '26.94'
$rng.Style = $origStyle
$ws.Range('E32').Value = '  -0.70%  '
$ws.Range('E33').Value = '  -3.21%  '
$ws.Range('E34').Value = '  +7.89%  '
$ws.Range('E35').Value = '  +2.41%  '
$ws.Range('E36').Value = '  +1.08%  '
$rng = $ws.Range('D37')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '3.41'
$rng.Style = $origStyle
$ws.Range('E37').Value = '  +12.16%  '
$rng = $ws.Range('D38')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '6.07'
$rng.Style = $origStyle
$ws.Range('E38').Value = '  -0.08%  '
$ws.Range('E39').Value = '  +0.14%  '
$rng = $ws.Range('D40')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '447.41'
$rng.Style = $origStyle
$ws.Range('E40').Value = '  +4.12%  '
$rng = $ws.Range('D41')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '8.67'
$rng.Style = $origStyle
$ws.Range('E41').Value = '  -1.45%  '
$rng = $ws.Range('D42')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '0.0371'
$rng.Style = $origStyle
$ws.Range('E42').Value = '  -1.07%  '
$ws.Range('D43').Value = '2.878.11'
$ws.Range('E43').Value = '  -2.71%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('E45').Value = '  -0.93%  '
$ws.Range('E46').Value = '  -0.34%  '
$rng = $ws.Range('D47')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '35.74'
$rng.Style = $origStyle
$ws.Range('E47').Value = '  +2.70%  '
$rng = $ws.Range('D49')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '123.30'
$rng.Style = $origStyle
$ws.Range('E49').Value = '  -2.14%  '
$ws.Range('E50').Value = '  -0.68%  '
$rng = $ws.Range('D51')
$origStyle = $rng.Style
$rng.NumberFormat = "@"
$rng.Value = '24.63'
$rng.Style = $origStyle
$ws.Range('E51').Value = '  -0.67%  '
